# Update sprint 4 velocity spreadsheet
# - Fill in the missing "Points Completed" value for Sprint 4 (row 5)
# - Add a new Sprint 5 row (row 6) with just the sprint number filled in
# - Leave the selection on A7, just below the newly added data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 4's "Points Completed" (C5) was left blank before; record it now.
$ws.Range("C5").Value = 28

# Start tracking Sprint 5: only the sprint number is known so far.
$ws.Range("A6").Value = 5

# Move the active selection to A7, ready for the next entry.
$ws.Range("A7").Select()
